$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V (home/away/odds/url) details between the given row pairs.
# A (Indice), B (pais), C (torneio), D (temporada) and E (data_partida) stay put;
# only the match details from column F through V trade places.
function Swap-Rows($r1, $r2) {
    $vals1 = @{}
    $vals2 = @{}
    for ($c = 6; $c -le 22; $c++) {
        $vals1[$c] = $ws.Cells.Item($r1, $c).Value2
        $vals2[$c] = $ws.Cells.Item($r2, $c).Value2
    }
    for ($c = 6; $c -le 22; $c++) {
        $ws.Cells.Item($r1, $c).Value = $vals2[$c]
        $ws.Cells.Item($r2, $c).Value = $vals1[$c]
    }
}

Swap-Rows 23 24
Swap-Rows 39 40
Swap-Rows 60 61
Swap-Rows 74 75
Swap-Rows 77 78

# --- Append two new match rows (80 and 81), copying row 79's formatting
# (bold/bordered index cell + datetime-formatted data_partida cell) down first.
$ws.Range("A79:V79").Copy()
$ws.Range("A80:V81").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "portugal"
$ws.Cells.Item(80, 3).Value = "liga-portugal-2"
$ws.Cells.Item(80, 4).Value = "2023-2024"
$ws.Cells.Item(80, 5).Value = 45235.625
$ws.Cells.Item(80, 6).Value = "Leiria"
$ws.Cells.Item(80, 7).Value = 4
$ws.Cells.Item(80, 8).Value = "Maritimo"
$ws.Cells.Item(80, 9).Value = 3
$ws.Cells.Item(80, 10).Value = 2.48
$ws.Cells.Item(80, 11).Value = "01/11/2023 16:12"
$ws.Cells.Item(80, 12).Value = 2.61
$ws.Cells.Item(80, 13).Value = "05/11/2023 14:53"
$ws.Cells.Item(80, 14).Value = 3.36
$ws.Cells.Item(80, 15).Value = "01/11/2023 16:12"
$ws.Cells.Item(80, 16).Value = 3.37
$ws.Cells.Item(80, 17).Value = "05/11/2023 14:53"
$ws.Cells.Item(80, 18).Value = 2.78
$ws.Cells.Item(80, 19).Value = "01/11/2023 16:12"
$ws.Cells.Item(80, 20).Value = 2.83
$ws.Cells.Item(80, 21).Value = "05/11/2023 14:53"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/leiria-maritimo/4CJfTUT4/"

$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "portugal"
$ws.Cells.Item(81, 3).Value = "liga-portugal-2"
$ws.Cells.Item(81, 4).Value = "2023-2024"
$ws.Cells.Item(81, 5).Value = 45235.79166666666
$ws.Cells.Item(81, 6).Value = "Tondela"
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = "Torreense"
$ws.Cells.Item(81, 9).Value = 2
$ws.Cells.Item(81, 10).Value = 2.56
$ws.Cells.Item(81, 11).Value = "01/11/2023 16:12"
$ws.Cells.Item(81, 12).Value = 2.3
$ws.Cells.Item(81, 13).Value = "05/11/2023 18:52"
$ws.Cells.Item(81, 14).Value = 3.22
$ws.Cells.Item(81, 15).Value = "01/11/2023 16:12"
$ws.Cells.Item(81, 16).Value = 3.2
$ws.Cells.Item(81, 17).Value = "05/11/2023 18:55"
$ws.Cells.Item(81, 18).Value = 2.79
$ws.Cells.Item(81, 19).Value = "01/11/2023 16:12"
$ws.Cells.Item(81, 20).Value = 3.49
$ws.Cells.Item(81, 21).Value = "05/11/2023 18:52"
$ws.Cells.Item(81, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/tondela-torreense/balfVjah/"
